$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.907.06"
$ws.Range("E2").Value = "  +3.03%  "
$ws.Range("D3").Value = "2.667.44"
$ws.Range("E3").Value = "  +3.16%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.40"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.45"
$ws.Range("E6").Value = "  +5.64%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("D9").Value = "2.663.02"
$ws.Range("E9").Value = "  +3.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("E10").Value = "  +15.19%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.09"
$ws.Range("E14").Value = "  +4.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000190"
$ws.Range("E15").Value = "  +7.57%  "
$ws.Range("D16").Value = "3.143.97"
$ws.Range("E16").Value = "  +2.68%  "
$ws.Range("D17").Value = "68.746.16"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "2.658.98"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.42"
$ws.Range("E19").Value = "  +4.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "367.03"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.46"
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  +5.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.62"
$ws.Range("E25").Value = "  +7.31%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.17"
$ws.Range("E27").Value = "  +3.64%  "
$ws.Range("E28").Value = "  +9.90%  "
$ws.Range("D29").Value = "2.779.37"
$ws.Range("E29").Value = "  +2.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "585.49"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +5.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.02"
$ws.Range("E33").Value = "  +6.29%  "
$ws.Range("E34").Value = "  +4.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.131"
$ws.Range("E35").Value = "  +6.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +5.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.44"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("E39").Value = "  +7.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.38"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("E41").Value = "  +5.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.369"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.68"
$ws.Range("E43").Value = "  +8.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.75"
$ws.Range("E44").Value = "  +5.84%  "
$ws.Range("D45").Value = "0.0₆0323"
$ws.Range("E45").Value = "  +15.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.76"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.80"
$ws.Range("E48").Value = "  +2.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.75"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.73"
$ws.Range("E50").Value = "  +3.80%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.17"
$ws.Range("E51").Value = "  +5.08%  "
